# "add property for block"
# Adds 12 new columns (Grass1-5, Crack1-5, Treasure1-2) to the XML-mapped
# table on Sheet1, extending it from A1:J53 to A1:V53, and updates the
# sheet view / column widths to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$lo = $ws.ListObjects.Item(1)

# Grow the table/autofilter range first (columns come in as generic
# "ColumnN" placeholders at this point).
$lo.Resize($ws.Range("A1:V53"))

# Now stamp the real header text in - this is what the table column
# names get derived from when the workbook is saved.
$headers = @("Grass1","Grass2","Grass3","Grass4","Grass5","Crack1","Crack2","Crack3","Crack4","Crack5","Treasure1","Treasure2")
$startCol = 11
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# Column widths for the newly added columns (best-effort; COM rounds
# widths to whole pixels internally so exact fractional widths from a
# hand-authored workbook can't always be reproduced bit-for-bit).
# Column J already had the right width/bestFit before this edit, so it
# is left untouched - only K is new.
$ws.Range("K1").EntireColumn.ColumnWidth = 13.2857142857
$ws.Range("L1:M1").EntireColumn.ColumnWidth = 14.5714285714
$ws.Range("N1").EntireColumn.ColumnWidth = 15
$ws.Range("O1").EntireColumn.ColumnWidth = 13.1428571429
$ws.Range("P1").EntireColumn.ColumnWidth = 16.4285714286
$ws.Range("Q1:S1").EntireColumn.ColumnWidth = 9.8571428571
$ws.Range("U1").EntireColumn.ColumnWidth = 15.7142857143
$ws.Range("V1").EntireColumn.ColumnWidth = 14.7142857143

# Move the view over to the new columns and select U7, matching the
# edited sheetView/selection in the workbook.
$ws.Range("I1").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("U7").Select() | Out-Null
